$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first four data rows (2006年-2009年), rows 2-5,
# so that the remaining rows (2010年-2016年) shift up to rows 2-8.
$ws.Range("A2:I5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
